$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data snapshot
$ws.Range("D2").Value = '29.732.25'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.888.48'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''0.7918'
$ws.Range("E5").Value = '  -1.90%  '
$ws.Range("D6").Value = '''241.82'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("D9").Value = '''25.47'
$ws.Range("E9").Value = '  -3.58%  '
$ws.Range("D10").Value = '''0.07022'
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").Value = '''0.08044'
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '''0.7647'
$ws.Range("E12").Value = '  +2.76%  '
$ws.Range("D13").Value = '1.891.31'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").Value = '''5.303'
$ws.Range("E14").Value = '  +2.26%  '
$ws.Range("D15").Value = '''92.12'
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").Value = '29.755.65'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = '''5.922'
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").Value = '''242.54'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").Value = '''0.000007682'
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("D21").Value = '''8.184'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '2.141.01'
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = '''0.1621'
$ws.Range("E25").Value = '  +8.25%  '
$ws.Range("D26").Value = '''9.296'
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").Value = '''163.67'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").Value = '''18.62'
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("D29").Value = '''2.052'
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("D30").Value = '''1.377'
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("D31").Value = '''1.532'
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("D32").Value = '''4.425'
$ws.Range("E32").Value = '  +3.03%  '
$ws.Range("D33").Value = '''0.05687'
$ws.Range("E33").Value = '  +2.67%  '
$ws.Range("D34").Value = '''4.073'
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = '''1.262'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '''0.7364'
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("D38").Value = '''2.712'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '''0.01913'
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").Value = '''2.772'
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("D41").Value = '''0.4392'
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = '''72.24'
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").Value = '''5.837'
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = '''0.8390'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '1.030.99'
$ws.Range("E46").Value = '  +5.02%  '
$ws.Range("E47").Value = '  +1.13%  '
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").Value = '''9.869'
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").Value = '''7.444'
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("D51").Value = '2.047.80'
$ws.Range("E51").Value = '  -0.48%  '
